# "Generate Report for Handoff"
# The b.md file has now been handed off for localization, so its status
# rows across the Overview / zh-cn / de-de sheets move from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# handoff file names / timestamps and a new warning about a stale handback.

$wb = $excel.ActiveWorkbook

$warning = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7906ce2040072942416d6d195d3d755cf595e9ac/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76a8088521c892143bffb1bc7acd631568063747/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md entry ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-19 02:35:39"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
# F3 goes from "True" to the literal text "False" -- assigning the bare
# string would get auto-typed to a real Excel boolean, so copy the
# existing text "False" cell (O2) which keeps it a text/shared-string cell.
$zh.Range("O2").Copy($zh.Range("F3"))
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-19 02:35:34"
$zh.Range("P3").Value = $warning
$zh.Columns.Item(16).ColumnWidth = 39.15

# ---- de-de sheet: row 3 is the b.md entry ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("O2").Copy($de.Range("F3"))
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-19 02:35:39"
$de.Range("P3").Value = $warning
$de.Columns.Item(16).ColumnWidth = 39.15
